$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asignaturas")
$ws.Range("A1").Value = "Test"
